$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange / UpDown values for existing row 6
$ws.Range("X6").Value = -0.45999100000000226
$ws.Range("Y6").Value = "Down"

# Append a brand new data row (row 7) with a fresh scan result.
# Copy row 6's formatting down to row 7 first so date/percent styles match.
$ws.Range("A6:W6").Copy()
$ws.Range("A7:W7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A7").Value = 42648.888796296298
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = "Buy"
$ws.Range("D7").Value = 36
$ws.Range("E7").Value = 18541
$ws.Range("F7").Value = 3296
$ws.Range("G7").Value = 62
$ws.Range("H7").Value = 32
$ws.Range("I7").Value = 80
$ws.Range("J7").Value = 19
$ws.Range("K7").Value = 52040
$ws.Range("L7").Value = 391
$ws.Range("M7").Value = 206
$ws.Range("N7").Value = 109
$ws.Range("O7").Value = 26
$ws.Range("P7").Value = "Noun"
$ws.Range("Q7").Value = 58.594837935340642
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.1167
$ws.Range("T7").Value = 0.0079000000000000008
$ws.Range("U7").Value = 5.99
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = 0
